$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Credentials")

$ws.Range("A9").Value = "Function_Change_Password"

$ws.Range("A10").Value = "username"
$ws.Range("B10").Value = "password"

$ws.Range("A11").Value = "jkl"
$ws.Range("B11").Value = "jkl"

$ws.Range("C10").Value = "newpassword"

$ws.Range("C11").Value = "jkl1"
